$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 08:17:53"
$ws.Range("E3").Value = "2026-02-06 08:17:56"
$c = $ws.Range("H3")
$c.NumberFormat = "@"
$c.Value = "73%"
$c.NumberFormat = "general"
$ws.Range("K3").Value = "0.0 MJ/m2"
$ws.Range("N3").Value = "-4.5 °C 7:34 TU"
$ws.Range("O3").Value = "-2.5 °C"
$ws.Range("E4").Value = "2026-02-06 08:17:58"
$ws.Range("J4").Value = "994.4 hPa"
$ws.Range("K4").Value = "0.1 MJ/m2"
$ws.Range("E5").Value = "2026-02-06 08:18:01"
$ws.Range("J5").Value = "994.8 hPa"
$ws.Range("K5").Value = "0.1 MJ/m2"
$ws.Range("O5").Value = "7.6 °C"
$ws.Range("E6").Value = "2026-02-06 08:18:03"
$ws.Range("J6").Value = "995.9 hPa"
$ws.Range("K6").Value = "0.2 MJ/m2"
$ws.Range("E7").Value = "2026-02-06 08:18:06"
$ws.Range("J7").Value = "995.7 hPa"
$ws.Range("K7").Value = "0.2 MJ/m2"
$ws.Range("E8").Value = "2026-02-06 08:18:08"
$ws.Range("K8").Value = "0.2 MJ/m2"
$ws.Range("E9").Value = "2026-02-06 08:18:11"
$ws.Range("O9").Value = "1.5 °C"
$ws.Range("E10").Value = "2026-02-06 08:18:13"
$ws.Range("E11").Value = "2026-02-06 08:18:15"
$c = $ws.Range("H11")
$c.NumberFormat = "@"
$c.Value = "86%"
$c.NumberFormat = "general"
$ws.Range("J11").Value = "997.0 hPa"
$ws.Range("O11").Value = "3.4 °C"
$ws.Range("E12").Value = "2026-02-06 08:18:18"
$c = $ws.Range("H12")
$c.NumberFormat = "@"
$c.Value = "63%"
$c.NumberFormat = "general"
$ws.Range("K12").Value = "0.1 MJ/m2"
$ws.Range("E13").Value = "2026-02-06 08:18:21"
$c = $ws.Range("H13")
$c.NumberFormat = "@"
$c.Value = "92%"
$c.NumberFormat = "general"
$ws.Range("E14").Value = "2026-02-06 08:18:23"
$ws.Range("I14").Value = "0.5 mm"
$ws.Range("E15").Value = "2026-02-06 08:18:26"
$ws.Range("J15").Value = "995.0 hPa"
$ws.Range("K15").Value = "0.2 MJ/m2"
$ws.Range("E16").Value = "2026-02-06 08:18:28"
$c = $ws.Range("H16")
$c.NumberFormat = "@"
$c.Value = "94%"
$c.NumberFormat = "general"
$ws.Range("N16").Value = "3.1 °C 7:43 TU"
$ws.Range("O16").Value = "4.0 °C"
$ws.Range("E17").Value = "2026-02-06 08:18:31"
$ws.Range("J17").Value = "997.8 hPa"
$ws.Range("K17").Value = "0.1 MJ/m2"
$ws.Range("E18").Value = "2026-02-06 08:18:33"
$ws.Range("K18").Value = "0.1 MJ/m2"
$ws.Range("N18").Value = "-5.6 °C 7:38 TU"
$ws.Range("O18").Value = "-5.0 °C"
$ws.Range("E19").Value = "2026-02-06 08:18:35"
$c = $ws.Range("H19")
$c.NumberFormat = "@"
$c.Value = "95%"
$c.NumberFormat = "general"
$ws.Range("J19").Value = "998.0 hPa"
$ws.Range("K19").Value = "0.1 MJ/m2"
$ws.Range("O19").Value = "6.5 °C"
$ws.Range("E20").Value = "2026-02-06 08:18:38"
$c = $ws.Range("H20")
$c.NumberFormat = "@"
$c.Value = "73%"
$c.NumberFormat = "general"
$ws.Range("K20").Value = "0.2 MJ/m2"
$ws.Range("O20").Value = "-2.7 °C"
$ws.Range("E21").Value = "2026-02-06 08:18:41"
$c = $ws.Range("H21")
$c.NumberFormat = "@"
$c.Value = "89%"
$c.NumberFormat = "general"
$ws.Range("J21").Value = "996.0 hPa"
$ws.Range("K21").Value = "0.1 MJ/m2"
$ws.Range("N21").Value = "1.7 °C 7:41 TU"
$ws.Range("O21").Value = "4.0 °C"
$ws.Range("E22").Value = "2026-02-06 08:18:43"
$c = $ws.Range("H22")
$c.NumberFormat = "@"
$c.Value = "88%"
$c.NumberFormat = "general"
$ws.Range("K22").Value = "0.2 MJ/m2"
$ws.Range("E23").Value = "2026-02-06 08:18:45"
$ws.Range("J23").Value = "995.0 hPa"
$ws.Range("K23").Value = "0.1 MJ/m2"
$ws.Range("E24").Value = "2026-02-06 08:18:48"
$ws.Range("J24").Value = "994.0 hPa"
$ws.Range("K24").Value = "0.1 MJ/m2"
$ws.Range("M24").Value = "12.9 °C 7:48 TU"
$ws.Range("E25").Value = "2026-02-06 08:18:50"
$ws.Range("J25").Value = "997.1 hPa"
$ws.Range("K25").Value = "0.1 MJ/m2"
$ws.Range("E26").Value = "2026-02-06 08:18:53"
$ws.Range("O26").Value = "-1.7 °C"
$ws.Range("E27").Value = "2026-02-06 08:18:55"
$ws.Range("J27").Value = "994.8 hPa"
$ws.Range("K27").Value = "0.2 MJ/m2"
$ws.Range("E28").Value = "2026-02-06 08:18:58"
$ws.Range("J28").Value = "998.3 hPa"
$ws.Range("O28").Value = "1.6 °C"
$ws.Range("E29").Value = "2026-02-06 08:19:00"
$c = $ws.Range("H29")
$c.NumberFormat = "@"
$c.Value = "71%"
$c.NumberFormat = "general"
$ws.Range("K29").Value = "0.1 MJ/m2"
$ws.Range("O29").Value = "9.7 °C"
$ws.Range("E30").Value = "2026-02-06 08:19:03"
$ws.Range("K30").Value = "0.2 MJ/m2"
$ws.Range("O30").Value = "-4.0 °C"
$ws.Range("E31").Value = "2026-02-06 08:19:05"
$ws.Range("E32").Value = "2026-02-06 08:19:07"
$ws.Range("J32").Value = "996.2 hPa"
$ws.Range("K32").Value = "0.1 MJ/m2"
$ws.Range("E33").Value = "2026-02-06 08:19:10"
$c = $ws.Range("H33")
$c.NumberFormat = "@"
$c.Value = "97%"
$c.NumberFormat = "general"
$ws.Range("E34").Value = "2026-02-06 08:19:12"
$ws.Range("K34").Value = "0.1 MJ/m2"
$ws.Range("O34").Value = "6.0 °C"
$ws.Range("E35").Value = "2026-02-06 08:19:14"
$ws.Range("K35").Value = "0.1 MJ/m2"
$ws.Range("N35").Value = "-3.4 °C 7:59 TU"
$ws.Range("E36").Value = "2026-02-06 08:19:17"
$c = $ws.Range("H36")
$c.NumberFormat = "@"
$c.Value = "71%"
$c.NumberFormat = "general"
$ws.Range("J36").Value = "997.7 hPa"
$ws.Range("K36").Value = "0.0 MJ/m2"
$ws.Range("O36").Value = "10.4 °C"
